{"js": "// Word JS API (Office.js) edit script.\n// Body is the content of `async (context) => { ... }`.\n\nconst body = context.document.body;\n\n// 1) Insert \"revised and resubmitted \" right before \"manuscript titled\"\n//    in the sentence \"Please find our manuscript titled ...\".\nlet hits = body.search(\"manuscript titled \\\"Colonial Architecture\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"revised and resubmitted \", \"Before\");\n}\nawait context.sync();\n\n// 2) Append a new sentence after \"...Journal of Experimental Biology.\"\nhits = body.search(\"to be considered for publication as a Research Article in the Journal of Experimental Biology.\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\" We believe we have now addressed the reviewers\\u2019 concerns.\", \"After\");\n}\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# 1) Insert \"revised and resubmitted \" right before \"manuscript titled\"\n#    in the sentence \"Please find our manuscript titled ...\".\n$rng1 = $d.Content\n$rng1.Find.Execute(\"manuscript titled\") | Out-Null\n$rng1.Collapse(1)  # wdCollapseStart\n$rng1.InsertBefore(\"revised and resubmitted \")\n\n# 2) Append a new sentence after \"...Journal of Experimental Biology.\"\n$rng2 = $d.Content\n$rng2.Find.Execute(\"Journal of Experimental Biology.\") | Out-Null\n$rng2.Collapse(0)  # wdCollapseEnd\n$rng2.InsertAfter(\" We believe we have now addressed the reviewers\" + [char]0x2019 + \" concerns.\")\n"}
